$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F98").Value = 13886962
$ws.Range("G98").Value = 11076398
$ws.Range("H98").Value = 904872
$ws.Range("H132").Value = 935377
$ws.Range("M132").Value = 14066924
$ws.Range("M133").Value = 14469376
$ws.Range("M134").Value = 14605011
$ws.Range("M135").Value = 14690448
$ws.Range("M136").Value = 14866907
$ws.Range("M137").Value = 14861425
$ws.Range("M138").Value = 14670817
$ws.Range("M139").Value = 14467042
$ws.Range("D140").Value = 823668
$ws.Range("M140").Value = 14131822
$ws.Range("N140").Value = 8935504
$ws.Range("M141").Value = 13894119
$ws.Range("N141").Value = 8782019
$ws.Range("M142").Value = 13750025
$ws.Range("N142").Value = 8776287
$ws.Range("M143").Value = 13621894
$ws.Range("N143").Value = 8836729
$ws.Range("N144").Value = 9042537
$ws.Range("N145").Value = 9395157
$ws.Range("N146").Value = 9716021
$ws.Range("N147").Value = 9916852
$ws.Range("N148").Value = 10176485
$ws.Range("N149").Value = 10402066
$ws.Range("N150").Value = 10586466
$ws.Range("N151").Value = 10834379
$ws.Range("D186").Value = 590429
$ws.Range("G186").Value = 9817905
$ws.Range("H186").Value = 1209134
$ws.Range("M186").Value = 13889924
$ws.Range("N186").Value = 8211234
$ws.Range("M187").Value = 13966157
$ws.Range("N187").Value = 8231225
$ws.Range("M188").Value = 14177390
$ws.Range("N188").Value = 8311011
$ws.Range("M189").Value = 14260698
$ws.Range("N189").Value = 8315371
$ws.Range("M190").Value = 14423025
$ws.Range("N190").Value = 8370952
$ws.Range("G191").Value = 8670408
$ws.Range("H191").Value = 1073055
$ws.Range("M191").Value = 14528177
$ws.Range("N191").Value = 8368471
$ws.Range("M192").Value = 14532951
$ws.Range("N192").Value = 8290045
$ws.Range("M193").Value = 14761475
$ws.Range("N193").Value = 8365137
$ws.Range("F194").Value = 14729052
$ws.Range("M194").Value = 14844156
$ws.Range("N194").Value = 8388024
$ws.Range("M195").Value = 14889086
$ws.Range("N195").Value = 8325496
$ws.Range("F196").Value = 13097794
$ws.Range("M196").Value = 14967768
$ws.Range("N196").Value = 8319992
$ws.Range("M197").Value = 15089912
$ws.Range("N197").Value = 8341456
$ws.Range("M198").Value = 15141503
$ws.Range("E199").Value = 214623
$ws.Range("M199").Value = 15309047
$ws.Range("O199").Value = 2587924
$ws.Range("M200").Value = 15440657
$ws.Range("O200").Value = 2599976
$ws.Range("M201").Value = 15510431
$ws.Range("O201").Value = 2605553
$ws.Range("F202").Value = 10071201
$ws.Range("M202").Value = 15673741
$ws.Range("O202").Value = 2635467
$ws.Range("F203").Value = 9851451
$ws.Range("O203").Value = 2650965
$ws.Range("O204").Value = 2658787
$ws.Range("O205").Value = 2687336
$ws.Range("O206").Value = 2682390
$ws.Range("O207").Value = 2694741
$ws.Range("O208").Value = 2719324
$ws.Range("O209").Value = 2707766
$ws.Range("O210").Value = 2680331
$ws.Range("B232").Value = 1361143
$ws.Range("F232").Value = 11589066
$ws.Range("B619").Value = 2287300
$ws.Range("D619").Value = 658996
$ws.Range("E619").Value = 333946
$ws.Range("F619").Value = 13905402
$ws.Range("G619").Value = 12418359
$ws.Range("H619").Value = 3986445
$ws.Range("I619").Value = 344.56
$ws.Range("M619").Value = 97731036
$ws.Range("N619").Value = 11764860
$ws.Range("O619").Value = 13082589
$ws.Range("B620").Value = 2044501
$ws.Range("D620").Value = 590827
$ws.Range("E620").Value = 286265
$ws.Range("F620").Value = 12835900
$ws.Range("G620").Value = 10782693
$ws.Range("H620").Value = 3456209
$ws.Range("I620").Value = 343.67
$ws.Range("M620").Value = 82594132
$ws.Range("N620").Value = 10199574
$ws.Range("O620").Value = 12419030
$ws.Range("B621").Value = 1503165
$ws.Range("D621").Value = 460445
$ws.Range("E621").Value = 274102
$ws.Range("F621").Value = 12782568
$ws.Range("G621").Value = 10601583
$ws.Range("H621").Value = 3445899
$ws.Range("I621").Value = 347.81
$ws.Range("J621").Value = 37.15
$ws.Range("K621").Value = 24.03
$ws.Range("L621").Value = 333.26
$ws.Range("M621").Value = 69687456
$ws.Range("N621").Value = 9344695
$ws.Range("O621").Value = 11625423
$ws.Range("B622").Value = 1410070
$ws.Range("D622").Value = 478187
$ws.Range("E622").Value = 305616
$ws.Range("F622").Value = 10027144
$ws.Range("G622").Value = 8093354
$ws.Range("H622").Value = 2683996
$ws.Range("I622").Value = 352.04
$ws.Range("L622").Value = 340.1
$ws.Range("M622").Value = 58985962
$ws.Range("N622").Value = 8911175
$ws.Range("O622").Value = 8847975
$ws.Range("B623").Value = 1211349
$ws.Range("D623").Value = 353957
$ws.Range("E623").Value = 305976
$ws.Range("F623").Value = 8847721
$ws.Range("G623").Value = 6804426
$ws.Range("H623").Value = 2255490
$ws.Range("I623").Value = 350.24
$ws.Range("J623").Value = 34.87
$ws.Range("L623").Value = 344.45
$ws.Range("M623").Value = 52914178
$ws.Range("N623").Value = 8504011
$ws.Range("O623").Value = 6719970
$ws.Range("B624").Value = 1151226
$ws.Range("D624").Value = 328969
$ws.Range("E624").Value = 210991
$ws.Range("F624").Value = 7846184
$ws.Range("G624").Value = 6420069
$ws.Range("H624").Value = 2220710
$ws.Range("I624").Value = 364.68
$ws.Range("K624").Value = 18.54
$ws.Range("L624").Value = 347.33
$ws.Range("M624").Value = 48457408
$ws.Range("N624").Value = 8091830
$ws.Range("O624").Value = 5601364
$ws.Range("B625").Value = 1249531
$ws.Range("D625").Value = 414871
$ws.Range("E625").Value = 182182
$ws.Range("F625").Value = 7669547
$ws.Range("G625").Value = 6098061
$ws.Range("H625").Value = 2192522
$ws.Range("J625").Value = 40.31
$ws.Range("K625").Value = 17.8
$ws.Range("L625").Value = 349.78
$ws.Range("M625").Value = 44260461
$ws.Range("N625").Value = 7631632
$ws.Range("O625").Value = 4742095
$ws.Range("B626").Value = 1413668
$ws.Range("D626").Value = 529876
$ws.Range("E626").Value = 177774
$ws.Range("F626").Value = 9153312
$ws.Range("G626").Value = 7008471
$ws.Range("H626").Value = 2634696
$ws.Range("I626").Value = 390.24
$ws.Range("J626").Value = 40.18
$ws.Range("K626").Value = 17.28
$ws.Range("L626").Value = 352.67
$ws.Range("M626").Value = 41276789
$ws.Range("N626").Value = 7250663
$ws.Range("O626").Value = 4098103
$ws.Range("B627").Value = 905143
$ws.Range("D627").Value = 392661
$ws.Range("E627").Value = 139214
$ws.Range("F627").Value = 7685163
$ws.Range("G627").Value = 6243758
$ws.Range("H627").Value = 2415534
$ws.Range("I627").Value = 399.94
$ws.Range("J627").Value = 38.57
$ws.Range("K627").Value = 16.82
$ws.Range("L627").Value = 355.81
$ws.Range("M627").Value = 38624166
$ws.Range("N627").Value = 6877198
$ws.Range("O627").Value = 3604111
$ws.Range("B628").Value = 892242
$ws.Range("D628").Value = 370108
$ws.Range("E628").Value = 168482
$ws.Range("F628").Value = 7548038
$ws.Range("G628").Value = 6438790
$ws.Range("H628").Value = 2540724
$ws.Range("I628").Value = 407.97
$ws.Range("J628").Value = 35.5
$ws.Range("K628").Value = 16.51
$ws.Range("L628").Value = 359.77
$ws.Range("M628").Value = 35818429
$ws.Range("N628").Value = 6397324
$ws.Range("O628").Value = 3163825
$ws.Range("B629").Value = 877095
$ws.Range("D629").Value = 304709
$ws.Range("E629").Value = 125343
$ws.Range("F629").Value = 6041681
$ws.Range("G629").Value = 4790100
$ws.Range("H629").Value = 1886936
$ws.Range("I629").Value = 403.98
$ws.Range("J629").Value = 33.87
$ws.Range("K629").Value = 17.38
$ws.Range("L629").Value = 363.43
$ws.Range("M629").Value = 33636726
$ws.Range("N629").Value = 5635760
$ws.Range("O629").Value = 2880436
$ws.Range("B630").Value = 859355
$ws.Range("D630").Value = 313024
$ws.Range("E630").Value = 125577
$ws.Range("F630").Value = 5994282
$ws.Range("G630").Value = 4675798
$ws.Range("H630").Value = 1815114
$ws.Range("I630").Value = 397.15
$ws.Range("J630").Value = 32.57
$ws.Range("K630").Value = 17.39
$ws.Range("L630").Value = 367.89
$ws.Range("M630").Value = 31534275
$ws.Range("N630").Value = 5196630
$ws.Range("O630").Value = 2635468
$ws.Range("B631").Value = 937326
$ws.Range("D631").Value = 327883
$ws.Range("E631").Value = 106319
$ws.Range("F631").Value = 5524780
$ws.Range("G631").Value = 4383369
$ws.Range("H631").Value = 1692387
$ws.Range("I631").Value = 394.63
$ws.Range("J631").Value = 31.55
$ws.Range("K631").Value = 16.92
$ws.Range("L631").Value = 372.77
$ws.Range("M631").Value = 29240217
$ws.Range("N631").Value = 4865517
$ws.Range("O631").Value = 2407841
$ws.Range("R631").Value = 324999
$ws.Range("B632").Value = 1026178
$ws.Range("D632").Value = 367178
$ws.Range("E632").Value = 113875
$ws.Range("F632").Value = 6322433
$ws.Range("G632").Value = 4669006
$ws.Range("H632").Value = 1754749
$ws.Range("I632").Value = 384.03
$ws.Range("J632").Value = 30.83
$ws.Range("K632").Value = 16.42
$ws.Range("L632").Value = 377.44
$ws.Range("M632").Value = 27538757
$ws.Range("N632").Value = 4641868
$ws.Range("O632").Value = 2235451
$ws.Range("P632").Value = 5246
$ws.Range("Q632").Value = 14157
$ws.Range("R632").Value = 331601
$ws.Range("B633").Value = 885265
$ws.Range("D633").Value = 355160
$ws.Range("E633").Value = 120804
$ws.Range("F633").Value = 6364213
$ws.Range("G633").Value = 5290199
$ws.Range("H633").Value = 2005846
$ws.Range("I633").Value = 389.97
$ws.Range("J633").Value = 30.28
$ws.Range("K633").Value = 15.63
$ws.Range("L633").Value = 382.63
$ws.Range("M633").Value = 26098703
$ws.Range("N633").Value = 4536583
$ws.Range("O633").Value = 2082153
$ws.Range("P633").Value = 5247
$ws.Range("Q633").Value = 14198

# New cells in row 633
$ws.Range("C633").Value = 1399.5
$ws.Range("R633").Value = 330337

# New row 634
$ws.Range("A634").Value = 44834
$ws.Range("B634").Value = 717451
$ws.Range("D634").Value = 285950
$ws.Range("E634").Value = 106802
$ws.Range("F634").Value = 5706132
$ws.Range("G634").Value = 4186917
$ws.Range("H634").Value = 1652775
$ws.Range("I634").Value = 403.8
$ws.Range("J634").Value = 29.44
$ws.Range("K634").Value = 15.42
$ws.Range("L634").Value = 387.52
$ws.Range("M634").Value = 25067483
$ws.Range("N634").Value = 4344346
$ws.Range("O634").Value = 1883339
$ws.Range("P634").Value = 5231
$ws.Range("Q634").Value = 14187
